# Edit script: update "Estado de Cuenta" worksheet with new workers/periods data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 2 extra data rows so the table grows from 7 to 9 rows ---
# Original data rows: 16-22 (7 rows). New data rows: 16-24 (9 rows).
# Insert two new rows right before the old last row (22), pushing it (and the
# signature block below) down by two rows.
$ws.Rows("22:23").Insert()

# The newly inserted rows 22 and 23 are blank/unformatted; copy formatting
# (and content, which we'll overwrite) from row 21 (a normal-style data row)
$ws.Range("B21:J21").Copy($ws.Range("B22:J22"))
$ws.Range("B21:J21").Copy($ws.Range("B23:J23"))

# Row 24 (now holding what used to be row22's content/style - the special
# "last row" bottom-border style) already has the correct style after the
# insert shifted it down, so nothing further needed there.

# --- 2. Fill in the data table (rows 16-24) ---
# Worker CARLOS ANDRES BALDIRIS REINO (CC 1047393841) - periods 2505,2406,2405,2404,2403
$ws.Range("B16").Value2 = "CC"
$ws.Range("C16").Value2 = "1047393841"
$ws.Range("D16").Value2 = "CARLOS ANDRES BALDIRIS REINO"
$ws.Range("E16").Value2 = "2505"
$ws.Range("F16").Value2 = 2388
$ws.Range("G16").Value2 = 1790932

$ws.Range("B17").Value2 = "CC"
$ws.Range("C17").Value2 = "1047393841"
$ws.Range("D17").Value2 = "CARLOS ANDRES BALDIRIS REINO"
$ws.Range("E17").Value2 = "2406"
$ws.Range("F17").Value2 = 52000
$ws.Range("G17").Value2 = 1790932

$ws.Range("B18").Value2 = "CC"
$ws.Range("C18").Value2 = "1047393841"
$ws.Range("D18").Value2 = "CARLOS ANDRES BALDIRIS REINO"
$ws.Range("E18").Value2 = "2405"
$ws.Range("F18").Value2 = 52000
$ws.Range("G18").Value2 = 1790932

$ws.Range("B19").Value2 = "CC"
$ws.Range("C19").Value2 = "1047393841"
$ws.Range("D19").Value2 = "CARLOS ANDRES BALDIRIS REINO"
$ws.Range("E19").Value2 = "2404"
$ws.Range("F19").Value2 = 52000
$ws.Range("G19").Value2 = 1790932

$ws.Range("B20").Value2 = "CC"
$ws.Range("C20").Value2 = "1047393841"
$ws.Range("D20").Value2 = "CARLOS ANDRES BALDIRIS REINO"
$ws.Range("E20").Value2 = "2403"
$ws.Range("F20").Value2 = 34666
$ws.Range("G20").Value2 = 1790932

# Worker ERIKA PATRICIA MATURANA ROSENSTAND (CC 1047372612) - periods 2208,2207,2206
$ws.Range("B21").Value2 = "CC"
$ws.Range("C21").Value2 = "1047372612"
$ws.Range("D21").Value2 = "ERIKA PATRICIA MATURANA ROSENSTAND"
$ws.Range("E21").Value2 = "2208"
$ws.Range("F21").Value2 = 68000
$ws.Range("G21").Value2 = 1700000

$ws.Range("B22").Value2 = "CC"
$ws.Range("C22").Value2 = "1047372612"
$ws.Range("D22").Value2 = "ERIKA PATRICIA MATURANA ROSENSTAND"
$ws.Range("E22").Value2 = "2207"
$ws.Range("F22").Value2 = 68000
$ws.Range("G22").Value2 = 1700000

$ws.Range("B23").Value2 = "CC"
$ws.Range("C23").Value2 = "1047372612"
$ws.Range("D23").Value2 = "ERIKA PATRICIA MATURANA ROSENSTAND"
$ws.Range("E23").Value2 = "2206"
$ws.Range("F23").Value2 = 68000
$ws.Range("G23").Value2 = 1700000

# Worker CAROLINA MARIMON SIMARRA (CC 1001833348) - period 2507 - new worker, last row
$ws.Range("B24").Value2 = "CC"
$ws.Range("C24").Value2 = "1001833348"
$ws.Range("D24").Value2 = "CAROLINA MARIMON SIMARRA"
$ws.Range("E24").Value2 = "2507"
$ws.Range("F24").Value2 = 68000
$ws.Range("G24").Value2 = 1700000

# --- 3. Update the summary header values ---
$ws.Range("E11").Value2 = 465054      # VALOR MORA total
$ws.Range("C13").Value2 = 3           # Cant. Trabajadores (was 2)
$ws.Range("F13").Value2 = 9           # Cant. Periodos (was 7)
